# ===================================================================
# State Balancing - First workflows
# Adds a new "SC" sheet and populates the CO DR-0100 XML workflow sheet
# ===================================================================
$wb = $excel.ActiveWorkbook

# --- Populate the "CO DR-0100 XML" sheet (previously empty) ---
$ws8 = $wb.Worksheets.Item("CO DR-0100 XML")

$b8_0 = New-Object 'object[,]' 16,6
$b8_0[0,0] = 'Inputs that we need'
$b8_0[1,1] = 'CO DR-0100 XML tool.xlsm'
$b8_0[1,5] = 'XLSM file with macros that will help us to get state balancing'
$b8_0[2,1] = 'Detail worksheet'
$b8_0[2,5] = 'This detail file is downloaded from Taxsolver'
$b8_0[5,0] = 'Step by step'
$b8_0[6,1] = '1st part'
$b8_0[7,1] = 1
$b8_0[7,2] = 'Go to the detail tab'
$b8_0[8,1] = 2
$b8_0[8,2] = 'Go to DR 0100 XML form page'
$b8_0[9,1] = 3
$b8_0[9,2] = 'click on "Export Detail as Worksheet"'
$b8_0[10,2] = 3.1
$b8_0[10,3] = 'Save the file in some folder in the P drive'
$b8_0[11,2] = 3.2
$b8_0[11,3] = 'This will be just a temp file'
$b8_0[12,1] = 4
$b8_0[12,2] = 'Copy and paste the Detail worksheet file into the Tool file (in the Detail worksheet)'
$b8_0[14,1] = '2nd part'
$b8_0[15,1] = 1
$b8_0[15,2] = 'Go to Information Sheet in Taxsolver'
$ws8.Range("A1:F16").Value = $b8_0

$b8_1 = New-Object 'object[,]' 16,8
$b8_1[0,0] = 2
$b8_1[0,1] = 'Get the State Registration ID'
$b8_1[0,5] = 'State ID'
$b8_1[1,0] = 3
$b8_1[1,1] = 'Open browser'
$b8_1[1,5] = 'https://www.colorado.gov/revenueonline/'
$b8_1[2,0] = 4
$b8_1[2,1] = 'Click on Sales and Use Tax --- Find Sales and Use Tax rates'
$b8_1[2,7] = 'At the bottom of the page'
$b8_1[3,0] = 5
$b8_1[3,1] = 'Click on View Business Location rates'
$b8_1[4,0] = 6
$b8_1[4,1] = 'Type the State ID (2) into the Colorado Account Number field'
$b8_1[5,0] = 7
$b8_1[5,1] = 'Wait for the page to load and then, click on Export'
$b8_1[6,0] = 8
$b8_1[6,1] = 'Download the file'
$b8_1[7,1] = 8.1
$b8_1[7,2] = 'Check which browser to use and what configuration each user needs to do before running the bot'
$b8_1[8,1] = 8.1999999999999993
$b8_1[8,2] = 'You can save the file in any location in your local drive'
$b8_1[9,1] = 8.3000000000000007
$b8_1[9,2] = 'The output is a TXT file, which contains the colorado table info'
$b8_1[10,0] = 9
$b8_1[10,1] = 'Get the text from the TXT file and split by line and by semicolon (;)'
$b8_1[11,1] = 9.1
$b8_1[11,2] = 'You must have at the end 19 columns'
$b8_1[12,0] = 10
$b8_1[12,1] = 'Paste the result datatable into the Tool file, Website worksheet.'
$b8_1[14,0] = '3hd part - Use XLSM tool file'
$b8_1[15,0] = 1
$b8_1[15,1] = 'Click on Unhide Columns button'
$ws8.Range("B17:I32").Value = $b8_1

$b8_2 = New-Object 'object[,]' 6,12
$b8_2[0,0] = 2
$b8_2[0,1] = 'Click on Hide Columns button'
$b8_2[1,0] = 3
$b8_2[1,1] = 'In the Rate-Services Fees worksheet, check if the sum of each column is equals to 0'
$b8_2[2,1] = '3.1 We can focus on the row 60, which starts with "Level County"'
$b8_2[2,11] = 'Ask Jay'
$b8_2[3,1] = '3.1 If some of them are different to 0, the bot should mark this as a failure'
$b8_2[4,0] = 4
$b8_2[4,1] = 'In the Balancing Sheet, check if TAX Not balanced and EXCEPT Not balanced cells are equals to 0'
$b8_2[5,1] = '4.1 If some of them are different to 0, the bot should mark this as a failure'
$ws8.Range("B33:M38").Value = $b8_2

# Hyperlink for the Colorado Revenue Online website cell
$ws8.Hyperlinks.Add($ws8.Range("G18"), "https://www.colorado.gov/revenueonline/", [System.Reflection.Missing]::Value, "https://www.colorado.gov/revenueonline/", $ws8.Range("G18"))


# --- Add the new "SC" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws9.Name = "SC"

$b9_0 = New-Object 'object[,]' 12,6
$b9_0[0,0] = 'Inputs that we need'
$b9_0[1,1] = 'SC ST-389 tool (2).xlsx'
$b9_0[1,5] = 'XLSX file with macros and formulas that will help us to get state balancing'
$b9_0[2,1] = 'Detail worksheet'
$b9_0[2,5] = 'This detail file is downloaded from Taxsolver'
$b9_0[4,0] = 'Step by step'
$b9_0[5,1] = '1st part'
$b9_0[6,1] = 1
$b9_0[6,2] = 'Go to the detail tab'
$b9_0[7,1] = 2
$b9_0[7,2] = 'Go to SC ST 389 form page'
$b9_0[8,1] = 3
$b9_0[8,2] = 'click on "Export Detail as Worksheet"'
$b9_0[9,2] = 3.1
$b9_0[9,3] = 'Save the file in some folder in the P drive'
$b9_0[10,2] = 3.2
$b9_0[10,3] = 'This will be just a temp file'
$b9_0[11,1] = 4
$b9_0[11,2] = 'Copy and paste the Detail worksheet file into the Tool file (in the Detail worksheet)'
$ws9.Range("A1:F12").Value = $b9_0

# Selection bookkeeping on each touched sheet
$ws9.Range("H7").Select()

$wsLA = $wb.Worksheets.Item("LA")
$wsLA.Activate()
$wsLA.Range("B3").Select()

$ws8.Activate()
$ws8.Range("F25").Select()
